# Cotações atualizadas - 2025-11-03
# Adds a new row (60) with the quotation values for 2025-11-03 (Excel
# date serial 45964), extending the data range from A1:E59 to A1:E60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

# Column A: date value, using the same date/time number format as the
# rest of the "Data" column (e.g. A59).
$ws.Cells.Item($row, 1).Value = 45964
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

# Columns B-E: quotation values (stored as text, comma decimal separator,
# matching the format already used throughout the sheet).
$ws.Cells.Item($row, 2).Value = "15,4031"
$ws.Cells.Item($row, 3).Value = "11,1575"
$ws.Cells.Item($row, 4).Value = "15,4031"
$ws.Cells.Item($row, 5).Value = "15,4031"
